# Fix latency units in report sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: "Utility" -> "Utility (Percent)"
$ws.Range("O2").Value = "Utility (Percent)"

# Append unit suffix to Read Latency columns (I: min, J: max, K: average)
# Rows 3-12 and 14-23 use "msec", row 13 uses "usec"
for ($row = 3; $row -le 23; $row++) {
    if ($row -eq 13) {
        $unit = " usec"
    } else {
        $unit = " msec"
    }

    foreach ($col in @("I", "J", "K")) {
        $cell = $ws.Range("$col$row")
        $cell.Value = "$($cell.Value2)$unit"
    }
}
